$wb = $excel.ActiveWorkbook

# --- Sheet3 (BSM): row 107 value updates ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("H107").Value = 999999
$ws3.Range("I107").Value = 999999
$ws3.Range("K107").Value = 999999
$ws3.Range("M107").Value = -998079

# --- Sheet4 (CRP): clear H:N for rows 129-141 except 136 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("H129:N129").ClearContents()
$ws4.Range("H130:N130").ClearContents()
$ws4.Range("H131:N131").ClearContents()
$ws4.Range("H132:N132").ClearContents()
$ws4.Range("H133:N133").ClearContents()
$ws4.Range("H134:N134").ClearContents()
$ws4.Range("H135:N135").ClearContents()
$ws4.Range("H137:N137").ClearContents()
$ws4.Range("H138:N138").ClearContents()
$ws4.Range("H139:N139").ClearContents()
$ws4.Range("H140:N140").ClearContents()
$ws4.Range("H141:N141").ClearContents()

# --- Sheet6 (GSM): add H:N values for rows 125-141 ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("H125").Value = 0
$ws6.Range("I125").Value = 0
$ws6.Range("J125").Value = 0
$ws6.Range("K125").Value = 0
$ws6.Range("L125").Value = 0
$ws6.Range("H126").Value = 3500
$ws6.Range("I126").Value = 3500
$ws6.Range("J126").Value = 0
$ws6.Range("K126").Value = 10500
$ws6.Range("L126").Value = 0
$ws6.Range("M126").Value = -8030
$ws6.Range("H127").Value = 0
$ws6.Range("I127").Value = 0
$ws6.Range("J127").Value = 0
$ws6.Range("K127").Value = 0
$ws6.Range("L127").Value = 0
$ws6.Range("H128").Value = 80000
$ws6.Range("I128").Value = 0
$ws6.Range("J128").Value = 80000
$ws6.Range("K128").Value = 0
$ws6.Range("L128").Value = 80000
$ws6.Range("N128").Value = -89960
$ws6.Range("H129").Value = 0
$ws6.Range("I129").Value = 0
$ws6.Range("J129").Value = 0
$ws6.Range("K129").Value = 0
$ws6.Range("L129").Value = 0
$ws6.Range("H130").Value = 0
$ws6.Range("I130").Value = 0
$ws6.Range("J130").Value = 0
$ws6.Range("K130").Value = 0
$ws6.Range("L130").Value = 0
$ws6.Range("H131").Value = 0
$ws6.Range("I131").Value = 0
$ws6.Range("J131").Value = 0
$ws6.Range("K131").Value = 0
$ws6.Range("L131").Value = 0
$ws6.Range("H132").Value = 4079.2
$ws6.Range("I132").Value = 3466
$ws6.Range("J132").Value = 4999
$ws6.Range("K132").Value = 10398
$ws6.Range("L132").Value = 14997
$ws6.Range("M132").Value = -7868
$ws6.Range("N132").Value = -20057
$ws6.Range("H133").Value = 0
$ws6.Range("I133").Value = 0
$ws6.Range("J133").Value = 0
$ws6.Range("K133").Value = 0
$ws6.Range("L133").Value = 0
$ws6.Range("H134").Value = 0
$ws6.Range("I134").Value = 0
$ws6.Range("J134").Value = 0
$ws6.Range("K134").Value = 0
$ws6.Range("L134").Value = 0
$ws6.Range("H135").Value = 0
$ws6.Range("I135").Value = 0
$ws6.Range("J135").Value = 0
$ws6.Range("K135").Value = 0
$ws6.Range("L135").Value = 0
$ws6.Range("H136").Value = 0
$ws6.Range("I136").Value = 0
$ws6.Range("J136").Value = 0
$ws6.Range("K136").Value = 0
$ws6.Range("L136").Value = 0
$ws6.Range("H137").Value = 0
$ws6.Range("I137").Value = 0
$ws6.Range("J137").Value = 0
$ws6.Range("K137").Value = 0
$ws6.Range("L137").Value = 0
$ws6.Range("H138").Value = 105000
$ws6.Range("I138").Value = 0
$ws6.Range("J138").Value = 105000
$ws6.Range("K138").Value = 0
$ws6.Range("L138").Value = 105000
$ws6.Range("N138").Value = -115280
$ws6.Range("H139").Value = 50000
$ws6.Range("I139").Value = 0
$ws6.Range("J139").Value = 50000
$ws6.Range("K139").Value = 0
$ws6.Range("L139").Value = 50000
$ws6.Range("N139").Value = -60280
$ws6.Range("H140").Value = 0
$ws6.Range("I140").Value = 0
$ws6.Range("J140").Value = 0
$ws6.Range("K140").Value = 0
$ws6.Range("L140").Value = 0
$ws6.Range("H141").Value = 0
$ws6.Range("I141").Value = 0
$ws6.Range("J141").Value = 0
$ws6.Range("K141").Value = 0
$ws6.Range("L141").Value = 0

# --- Sheet7 (LTW): clear H:N for rows 124-141 except 126 ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("H124:N124").ClearContents()
$ws7.Range("H125:N125").ClearContents()
$ws7.Range("H127:N127").ClearContents()
$ws7.Range("H128:N128").ClearContents()
$ws7.Range("H129:N129").ClearContents()
$ws7.Range("H130:N130").ClearContents()
$ws7.Range("H131:N131").ClearContents()
$ws7.Range("H132:N132").ClearContents()
$ws7.Range("H133:N133").ClearContents()
$ws7.Range("H134:N134").ClearContents()
$ws7.Range("H135:N135").ClearContents()
$ws7.Range("H136:N136").ClearContents()
$ws7.Range("H137:N137").ClearContents()
$ws7.Range("H138:N138").ClearContents()
$ws7.Range("H139:N139").ClearContents()
$ws7.Range("H140:N140").ClearContents()
$ws7.Range("H141:N141").ClearContents()
